# Update workbook to reflect data pulled through 2022-08-07.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Rename the worksheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-08-07"

# Update the column header label (shared string) for the 2022 YTD column.
$ws.Range("I1").Value = "2022 (through 08-07)"

# August 2022 YTD carjackings: 31 -> 38
$ws.Range("I9").Value = 38

# Total row for the 2022 YTD column: 1001 -> 1008
$ws.Range("I14").Value = 1008
